# "Generate Report for Handback"
#
# The localization-status report tracks, per source file, the handoff /
# handback status of each target-locale translation. This run records that
# the handback for f2bb82f3-1d78-43f1-8acb-c8ab633f696a.md has completed for
# both target locales: the status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and the "Latest Handback DateTime" on
# each locale's detail sheet is stamped with the handback timestamp.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet: roll up the new status for both locale columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack

# --- zh-cn detail sheet ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusHandedBack
$zhcn.Range("G3").Value = "2016-02-17 04:25:15"

# --- de-de detail sheet ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusHandedBack
$dede.Range("G3").Value = "2016-02-17 04:25:33"
